$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column "Price" cells look numeric (e.g. "594.20", "63.888.27") but must
# stay plain text, matching the workbook's original inlineStr storage, so a
# leading apostrophe forces text entry (like a user typing it in Excel) and
# the style is then reset to Normal to drop the transient quote-prefix flag.

# Row 2
$cell = $ws.Range("D2")
$cell.Value = "'" + '63.824.10'
$cell.Style = "Normal"
$ws.Range("E2").Value = '  +0.08%  '

# Row 3
$cell = $ws.Range("D3")
$cell.Value = "'" + '2.620.68'
$cell.Style = "Normal"
$ws.Range("E3").Value = '  -0.14%  '

# Row 5
$cell = $ws.Range("D5")
$cell.Value = "'" + '594.20'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  -0.23%  '

# Row 6
$cell = $ws.Range("D6")
$cell.Value = "'" + '151.01'
$cell.Style = "Normal"
$ws.Range("E6").Value = '  +0.72%  '

# Row 7
$ws.Range("E7").Value = '  +0.02%  '

# Row 8
$cell = $ws.Range("D8")
$cell.Value = "'" + '0.587'
$cell.Style = "Normal"
$ws.Range("E8").Value = '  -0.19%  '

# Row 9
$ws.Range("E9").Value = '  +4.48%  '

# Row 10
$ws.Range("B10").Value = 'Cardano'
$ws.Range("C10").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$cell = $ws.Range("D10")
$cell.Value = "'" + '0.395'
$cell.Style = "Normal"
$ws.Range("E10").Value = '  +3.36%  '

# Row 11
$ws.Range("B11").Value = 'Toncoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$cell = $ws.Range("D11")
$cell.Value = "'" + '5.79'
$cell.Style = "Normal"
$ws.Range("E11").Value = '  +1.96%  '

# Row 12
$ws.Range("E12").Value = '  +1.07%  '

# Row 13
$cell = $ws.Range("D13")
$cell.Value = "'" + '27.91'
$cell.Style = "Normal"
$ws.Range("E13").Value = '  +0.77%  '

# Row 14
$cell = $ws.Range("D14")
$cell.Value = "'" + '3.091.34'
$cell.Style = "Normal"
$ws.Range("E14").Value = '  -0.14%  '

# Row 15
$ws.Range("B15").Value = 'ShibaInu'
$ws.Range("C15").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$cell = $ws.Range("D15")
$cell.Value = "'" + '0.0000170'
$cell.Style = "Normal"
$ws.Range("E15").Value = '  +13.69%  '

# Row 16
$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$cell = $ws.Range("D16")
$cell.Value = "'" + '63.743.13'
$cell.Style = "Normal"
$ws.Range("E16").Value = '  +0.23%  '

# Row 17
$cell = $ws.Range("D17")
$cell.Value = "'" + '2.616.83'
$cell.Style = "Normal"
$ws.Range("E17").Value = '  +0.11%  '

# Row 18
$ws.Range("E18").Value = '  -0.53%  '

# Row 19
$ws.Range("E19").Value = '  +2.98%  '

# Row 20
$cell = $ws.Range("D20")
$cell.Value = "'" + '348.29'
$cell.Style = "Normal"
$ws.Range("E20").Value = '  -0.12%  '

# Row 21
$ws.Range("E21").Value = '  +2.29%  '

# Row 22
$cell = $ws.Range("D22")
$cell.Value = "'" + '0.999'
$cell.Style = "Normal"
$ws.Range("E22").Value = '  +0.14%  '

# Row 23
$cell = $ws.Range("D23")
$cell.Value = "'" + '67.30'
$cell.Style = "Normal"
$ws.Range("E23").Value = '  +1.47%  '

# Row 24
$ws.Range("E24").Value = '  -2.96%  '

# Row 25
$ws.Range("E25").Value = '  +0.47%  '

# Row 26
$cell = $ws.Range("D26")
$cell.Value = "'" + '9.18'
$cell.Style = "Normal"
$ws.Range("E26").Value = '  -0.25%  '

# Row 27
$cell = $ws.Range("D27")
$cell.Value = "'" + '8.29'
$cell.Style = "Normal"
$ws.Range("E27").Value = '  +0.89%  '

# Row 28
$cell = $ws.Range("D28")
$cell.Value = "'" + '548.57'
$cell.Style = "Normal"
$ws.Range("E28").Value = '  -2.14%  '

# Row 29
$cell = $ws.Range("D29")
$cell.Value = "'" + '0.161'
$cell.Style = "Normal"
$ws.Range("E29").Value = '  -1.84%  '

# Row 30
$ws.Range("E30").Value = '  -0.25%  '

# Row 31
$cell = $ws.Range("D31")
$cell.Value = "'" + '0.0₃0905'
$cell.Style = "Normal"
$ws.Range("E31").Value = '  +7.42%  '

# Row 32
$cell = $ws.Range("D32")
$cell.Value = "'" + '2.07'
$cell.Style = "Normal"
$ws.Range("E32").Value = '  +1.68%  '

# Row 33
$ws.Range("E33").Value = '  +4.63%  '

# Row 34
$cell = $ws.Range("D34")
$cell.Value = "'" + '5.45'
$cell.Style = "Normal"
$ws.Range("E34").Value = '  +4.45%  '

# Row 35
$ws.Range("E35").Value = '  +0.66%  '

# Row 36
$cell = $ws.Range("D36")
$cell.Value = "'" + '0.419'
$cell.Style = "Normal"
$ws.Range("E36").Value = '  +2.62%  '

# Row 37
$cell = $ws.Range("D37")
$cell.Value = "'" + '164.50'
$cell.Style = "Normal"
$ws.Range("E37").Value = '  -2.52%  '

# Row 38
$cell = $ws.Range("D38")
$cell.Value = "'" + '19.95'
$cell.Style = "Normal"
$ws.Range("E38").Value = '  +3.23%  '

# Row 39
$ws.Range("B39").Value = 'Stacks'
$ws.Range("C39").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$cell = $ws.Range("D39")
$cell.Value = "'" + '1.97'
$cell.Style = "Normal"
$ws.Range("E39").Value = '  +1.76%  '

# Row 40
$ws.Range("B40").Value = 'FirstDigitalUSD'
$ws.Range("C40").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$cell = $ws.Range("D40")
$cell.Value = "'" + '0.999'
$cell.Style = "Normal"
$ws.Range("E40").Value = '  -0.02%  '

# Row 41
$ws.Range("E41").Value = '  +0.02%  '

# Row 42
$cell = $ws.Range("D42")
$cell.Value = "'" + '167.39'
$cell.Style = "Normal"
$ws.Range("E42").Value = '  -1.23%  '

# Row 43
$ws.Range("E43").Value = '  +4.54%  '

# Row 44
$cell = $ws.Range("D44")
$cell.Value = "'" + '23.20'
$cell.Style = "Normal"
$ws.Range("E44").Value = '  +8.88%  '

# Row 45
$cell = $ws.Range("D45")
$cell.Value = "'" + '2.21'
$cell.Style = "Normal"
$ws.Range("E45").Value = '  +11.61%  '

# Row 46
$ws.Range("E46").Value = '  -1.93%  '

# Row 47
$cell = $ws.Range("D47")
$cell.Value = "'" + '0.636'
$cell.Style = "Normal"
$ws.Range("E47").Value = '  +1.09%  '

# Row 48
$ws.Range("E48").Value = '  +1.52%  '

# Row 49
$ws.Range("E49").Value = '  +0.39%  '

# Row 50
$cell = $ws.Range("D50")
$cell.Value = "'" + '19.23'
$cell.Style = "Normal"
$ws.Range("E50").Value = '  +0.47%  '

# Row 51
$ws.Range("E51").Value = '  +17.93%  '
